$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "GRID1"
$ws.Range("B3").Value = "GRID1"
$ws.Range("B4").Value = "PORT7"
$ws.Range("B5").Value = "GRID1"
$ws.Range("B6").Value = "GRID1"
$ws.Range("B1").Value = "Current_Location"

$ws.Range("B2").Select()
